# Update the "Förändrad" (Changed) date column (C) for all data rows.
# Every value in C2:C78 moves from serial date 45192 (2023-09-23) to
# serial date 45202 (2023-10-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 78 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45192) {
        $cell.Value2 = 45202
    }
}
